$d = $word.ActiveDocument

# --- Update the "Tempo total do processo" timeline numbers ---
$d.Content.Find.Execute("Concluído (126 dias)", $false, $false, $false, $false, $false, $true, 1, $false, "Concluído (132 dias)", 2) | Out-Null
$d.Content.Find.Execute("Assinatura Contrato (135 dias)", $false, $false, $false, $false, $false, $true, 1, $false, "Assinatura Contrato (141 dias)", 2) | Out-Null
$d.Content.Find.Execute("Assinatura Contrato (43 dias)", $false, $false, $false, $false, $false, $true, 1, $false, "Assinatura Contrato (49 dias)", 2) | Out-Null
$d.Content.Find.Execute("Concluído (51 dias)", $false, $false, $false, $false, $false, $true, 1, $false, "Concluído (57 dias)", 2) | Out-Null
$d.Content.Find.Execute("Assinatura Contrato (78 dias)", $false, $false, $false, $false, $false, $true, 1, $false, "Assinatura Contrato (84 dias)", 2) | Out-Null
$d.Content.Find.Execute("Total de dias 853", $false, $false, $false, $false, $false, $true, 1, $false, "Total de dias 883", 2) | Out-Null

# --- Remove the dead PDM / "Valor Homologado" itemised list paragraph ---
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "PDM 1199 - Apito*") {
        $para.Range.Text = ""
        break
    }
}
